$d = $word.ActiveDocument

# 1. "Demonstrating utility" -> "demonstrating utility" (lowercase the D)
$d.Content.Find.Execute("Demonstrating utility", $true, $true, $false, $false, $false, $true, 1, $false, "demonstrating utility", 2)

# 2. "here introduced" -> "here-introduced" (add hyphen)
$d.Content.Find.Execute("here introduced", $true, $true, $false, $false, $false, $true, 1, $false, "here-introduced", 2)

# 3. "MlLib" -> "MLlib" (fix capitalization)
$d.Content.Find.Execute("MlLib", $true, $true, $false, $false, $false, $true, 1, $false, "MLlib", 2)
